# B6-PowerPoint.pptx edit
#
# 1) Three tables (slides 14, 15 and 16 - each holding exactly one table as
#    their first shape) switch from the deck's custom table style
#    {67F619A2-B602-402B-A6B4-33E8EC55AE31} to the built-in table style
#    {4398D058-C173-4857-A17D-24DE9F8DA235}.
#
# 2) The presentation's theme colour scheme (the only theme part the PowerPoint
#    object model exposes here - it backs the single slide master / "Integral"
#    design) is restored to the stock "Office" palette.

$p = $ppt.ActivePresentation

# --- 1) retarget the three table styles -----------------------------------
$tableSlides = 14, 15, 16
$newStyleId = "{4398D058-C173-4857-A17D-24DE9F8DA235}"

foreach ($idx in $tableSlides) {
    $slide = $p.Slides.Item($idx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newStyleId)
    }
}

# --- 2) restore the stock "Office" theme colour palette --------------------
function RGBval([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$master = $p.Designs.Item(1).SlideMaster
$colors = $master.Theme.ThemeColorScheme

$colors.Item(1).RGB = RGBval 0x00 0x00 0x00   # dk1
$colors.Item(2).RGB = RGBval 0xFF 0xFF 0xFF   # lt1
$colors.Item(3).RGB = RGBval 0x44 0x54 0x6A   # dk2
$colors.Item(4).RGB = RGBval 0xE7 0xE6 0xE6   # lt2
$colors.Item(5).RGB = RGBval 0x5B 0x9B 0xD5   # accent1
$colors.Item(6).RGB = RGBval 0xED 0x7D 0x31   # accent2
$colors.Item(7).RGB = RGBval 0xA5 0xA5 0xA5   # accent3
$colors.Item(8).RGB = RGBval 0xFF 0xC0 0x00   # accent4
$colors.Item(9).RGB = RGBval 0x44 0x72 0xC4   # accent5
$colors.Item(10).RGB = RGBval 0x70 0xAD 0x47  # accent6
$colors.Item(11).RGB = RGBval 0x05 0x63 0xC1  # hlink
$colors.Item(12).RGB = RGBval 0x95 0x4F 0x72  # folHlink
